# Refresh the coin price (col D) and 1h volume-change (col E) figures
# with the latest values from the data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.343.54"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "3.752.82"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'594.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "'169.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "3.751.03"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.98%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").Value = "'6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "'36.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "4.382.26"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "3.751.41"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "67.250.20"
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'7.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'10.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.57%  "
$ws.Range("D22").Value = "'466.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'0.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("E24").Value = "  -8.57%  "
$ws.Range("D25").Value = "'83.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").Value = "'12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("D31").Value = "3.900.55"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'30.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").Value = "3.715.00"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "'3.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").Value = "'0.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'0.311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("D47").Value = "'45.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").Value = "'396.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("E49").Value = "  -8.70%  "
$ws.Range("D50").Value = "'138.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("E51").Value = "  -2.41%  "
